# Set-up boilerplate for BaseClass benchmarking
# Populate the "HasClassGenerator" (column F) checkbox cells for the
# CmsBaseClasses rows (ValidationHelperGet, ModuleRegistration,
# WebPartBase, UiWebPartBase, UserControlBase, PageBase) with TRUE.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F20 (BH2500 - ValidationHelperGet) and F22 (BH3000 - ModuleRegistration)
# keep their current formatting, just flip the value to TRUE.
$ws.Range("F20").Value = $true
$ws.Range("F22").Value = $true

# F23:F26 (BH3500 WebPartBase, BH3501 UiWebPartBase, BH3502 UserControlBase,
# BH3503 PageBase) pick up the bordered "good" style already used by the
# BH4000 row (F27) and get flipped to TRUE as well.
$ws.Range("F27").Copy() | Out-Null
$ws.Range("F23:F26").PasteSpecial(-4122) | Out-Null
$ws.Range("F23:F26").Value = $true

$excel.CutCopyMode = 0

# Move the active selection to reflect where the edits were made.
$ws.Range("F20").Select() | Out-Null
